# Rename the old "Main_240414" sheet to "240414".
# Excel automatically rewrites any formula references to the sheet
# (e.g. in Merge_RKTM) to use the new, quoted name '240414'.
$wb = $excel.ActiveWorkbook
$mainOld = $wb.Worksheets.Item("Main_240414")
$mainOld.Name = "240414"

# Add a brand-new worksheet that will become "Main_240511" and place it
# immediately before the (renamed) "240414" sheet, i.e. as the first tab.
$newMain = $wb.Worksheets.Add($mainOld)
$newMain.Name = "Main_240511"

$data = @(
    @('Class+Node [(Identifier (Key)]', 'Class [Not chosen]', 'Node [Not chosen]', 'Required Mods [Not chosen]', 'English [Source string]', 'Korean (한국어) [Translation]'),
    @('HediffDef+SimpleFinger.label', 'HediffDef', 'SimpleFinger.label', $null, 'simple finger', '단순한 손가락 모형'),
    @('HediffDef+SimpleFinger.description', 'HediffDef', 'SimpleFinger.description', $null, 'An installed simple finger.', '이식된 단순한 손가락 모형.'),
    @('HediffDef+SimpleToe.label', 'HediffDef', 'SimpleToe.label', $null, 'simple toe', '단순한 발가락 모형'),
    @('HediffDef+SimpleToe.description', 'HediffDef', 'SimpleToe.description', $null, 'An installed simple toe.', '이식된 단순한 발가락 모형.'),
    @('HediffDef+ProstheticFinger.label', 'HediffDef', 'ProstheticFinger.label', $null, 'prosthetic finger', '간단한 인공 손가락'),
    @('HediffDef+ProstheticFinger.description', 'HediffDef', 'ProstheticFinger.description', $null, 'An installed prosthetic finger.', '이식된 간단한 인공 손가락.'),
    @('HediffDef+ProstheticToe.label', 'HediffDef', 'ProstheticToe.label', $null, 'prosthetic toe', '간단한 인공 발가락'),
    @('HediffDef+ProstheticToe.description', 'HediffDef', 'ProstheticToe.description', $null, 'An installed prosthetic toe.', '이식된 간단한 인공 발가락.'),
    @('HediffDef+BionicFinger.label', 'HediffDef', 'BionicFinger.label', $null, 'bionic finger', '생체공학 손가락'),
    @('HediffDef+BionicFinger.description', 'HediffDef', 'BionicFinger.description', $null, 'An installed bionic finger.', '이식된 생체공학 손가락.'),
    @('HediffDef+BionicToe.label', 'HediffDef', 'BionicToe.label', $null, 'bionic toe', '생체공학 발가락'),
    @('HediffDef+BionicToe.description', 'HediffDef', 'BionicToe.description', $null, 'An installed bionic toe.', '이식된 생체공학 발가락.'),
    @('RecipeDef+InstallSimpleFinger.label', 'RecipeDef', 'InstallSimpleFinger.label', $null, 'install simple finger', '단순한 손가락 모형 이식'),
    @('RecipeDef+InstallSimpleFinger.description', 'RecipeDef', 'InstallSimpleFinger.description', $null, 'Install a simple finger.', '단순한 손가락 모형을 이식합니다.'),
    @('RecipeDef+InstallSimpleFinger.jobString', 'RecipeDef', 'InstallSimpleFinger.jobString', $null, 'Installing simple finger.', '이식 수술 중'),
    @('RecipeDef+InstallSimpleToe.label', 'RecipeDef', 'InstallSimpleToe.label', $null, 'install simple toe', '단순한 발가락 모형 이식'),
    @('RecipeDef+InstallSimpleToe.description', 'RecipeDef', 'InstallSimpleToe.description', $null, 'Install a simple toe.', '단순한 발가락 모형을 이식합니다.'),
    @('RecipeDef+InstallSimpleToe.jobString', 'RecipeDef', 'InstallSimpleToe.jobString', $null, 'Installing simple toe.', '이식 수술 중'),
    @('RecipeDef+InstallProstheticFinger.label', 'RecipeDef', 'InstallProstheticFinger.label', $null, 'install prosthetic finger', '인공 손가락 이식'),
    @('RecipeDef+InstallProstheticFinger.description', 'RecipeDef', 'InstallProstheticFinger.description', $null, 'Install a prosthetic finger.', '인공 손가락을 이식합니다.'),
    @('RecipeDef+InstallProstheticFinger.jobString', 'RecipeDef', 'InstallProstheticFinger.jobString', $null, 'Installing prosthetic finger.', '이식 수술 중'),
    @('RecipeDef+InstallProstheticToe.label', 'RecipeDef', 'InstallProstheticToe.label', $null, 'install prosthetic toe', '인공 발가락 이식'),
    @('RecipeDef+InstallProstheticToe.description', 'RecipeDef', 'InstallProstheticToe.description', $null, 'Install a prosthetic toe.', '인공 발가락을 이식합니다.'),
    @('RecipeDef+InstallProstheticToe.jobString', 'RecipeDef', 'InstallProstheticToe.jobString', $null, 'Installing prosthetic toe.', '이식 수술 중'),
    @('RecipeDef+InstallBionicFinger.label', 'RecipeDef', 'InstallBionicFinger.label', $null, 'install bionic finger', '생체공학 손가락 이식'),
    @('RecipeDef+InstallBionicFinger.description', 'RecipeDef', 'InstallBionicFinger.description', $null, 'Install a bionic finger.', '생체공학 손가락을 이식합니다.'),
    @('RecipeDef+InstallBionicFinger.jobString', 'RecipeDef', 'InstallBionicFinger.jobString', $null, 'Installing bionic finger.', '이식 수술 중'),
    @('RecipeDef+InstallBionicToe.label', 'RecipeDef', 'InstallBionicToe.label', $null, 'install bionic toe', '생체공학 발가락 이식'),
    @('RecipeDef+InstallBionicToe.description', 'RecipeDef', 'InstallBionicToe.description', $null, 'Install a bionic toe.', '생체공학 발가락을 이식합니다.'),
    @('RecipeDef+InstallBionicToe.jobString', 'RecipeDef', 'InstallBionicToe.jobString', $null, 'Installing bionic toe.', '이식 수술 중'),
    @('ThingDef+SimpleFinger.label', 'ThingDef', 'SimpleFinger.label', $null, 'simple finger replacement', '단순한 손가락 모형'),
    @('ThingDef+SimpleFinger.description', 'ThingDef', 'SimpleFinger.description', $null, 'A simple solid finger-like rod, makes gripping slightly easier.', '막대기같은 손가락 모형으로 쥐는 행동을 좀더 편하게 만들어줍니다.'),
    @('ThingDef+SimpleToe.label', 'ThingDef', 'SimpleToe.label', $null, 'simple toe replacement', '단순한 발가락 모형'),
    @('ThingDef+SimpleToe.description', 'ThingDef', 'SimpleToe.description', $null, 'A simple solid toe-like rod, makes walking slightly easier.', '막대기같은 발가락 모형으로 걷는것을 좀더 편하게 만들어줍니다.'),
    @('ThingDef+ProstheticFinger.label', 'ThingDef', 'ProstheticFinger.label', $null, 'prosthetic finger', '간단한 인공 손가락'),
    @('ThingDef+ProstheticFinger.description', 'ThingDef', 'ProstheticFinger.description', $null, 'A simple prosthetic finger in case of losing one of them. Slightly distinguishable from the real one.', '간단하지만 튼튼한 보철 손가락입니다.'),
    @('ThingDef+ProstheticToe.label', 'ThingDef', 'ProstheticToe.label', $null, 'prosthetic toe', '간단한 인공 발가락'),
    @('ThingDef+ProstheticToe.description', 'ThingDef', 'ProstheticToe.description', $null, 'The simplest prosthetic for a human toe. It feels occasionally inconvenient.', '몇개 잘린 발가락을 보충해줄 수 있는 간단한 보철 발가락입니다.'),
    @('ThingDef+BionicFinger.label', 'ThingDef', 'BionicFinger.label', $null, 'bionic finger', '생체공학 손가락'),
    @('ThingDef+BionicFinger.description', 'ThingDef', 'BionicFinger.description', $null, 'A subtle replacement for a human finger. Those are better than biological ones in almost every way.', '멋진 생체공학 손가락!'),
    @('ThingDef+BionicToe.label', 'ThingDef', 'BionicToe.label', $null, 'bionic toe', '생체공학 발가락'),
    @('ThingDef+BionicToe.description', 'ThingDef', 'BionicToe.description', $null, 'A rather good bionic toe replacement. It transmits electrical impulses to a biogel neural interface that communicates directly with the user''s nervous system.', '멋진 생체공학 발가락!')
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $v = $rowVals[$c]
        if ($null -ne $v) {
            $newMain.Cells.Item($r + 1, $c + 1).Value = $v
        }
    }
}

$newMain.Columns.Item(1).ColumnWidth = 44.7265625
$newMain.Columns.Item(2).ColumnWidth = 19.1796875
$newMain.Columns.Item(3).ColumnWidth = 33.453125
$newMain.Columns.Item(4).ColumnWidth = 29.26953125
$newMain.Columns.Item(5).ColumnWidth = 44.453125
$newMain.Columns.Item(6).ColumnWidth = 68.6328125

# Make the new sheet the active / selected tab, matching the authored
# workbook where "Main_240511" is shown selected (tabSelected) and the
# previous main sheet no longer is.
$newMain.Activate()
$newMain.Range("F44").Select()
